# "add tree, add lc.94,144,145"
#
# Sheet2!B9 is the underlying data cell that feeds (via formulas) into
# Sheet1!B9, Sheet1!C10 and Sheet1!D17. Its value moves from "8" to the
# new leetcode-style entry "94,144,145". Sheet1!B9 (today's row) also
# gets the "highlighted" cell style that is already used elsewhere in
# the column (e.g. B2), and the active sheet/selection bookmarks move
# from Sheet2 (C10 selected) to Sheet1 (B9 selected).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# 1. Update the source data cell. Sheet1!B9/C10/D17 all reference
#    Sheet2!B9 via formulas, so they recalculate automatically.
$ws2.Range("B9").Value = "94,144,145"

# 2. Give Sheet1!B9 the highlighted ("today") look by copying the
#    formatting only (not the value/formula) from a cell that already
#    carries that exact style.
$ws1.Range("B2").Copy()
$ws1.Range("B9").PasteSpecial(-4122)

# 3. Move the selection on Sheet2 to B10 (leaving Sheet2 not the active
#    tab), then activate Sheet1 and select B9 there so Sheet1 ends up
#    as the active/visible sheet.
$ws2.Activate()
$ws2.Range("B10").Select()
$ws1.Activate()
$ws1.Range("B9").Select()
